# Auto-generated edit script: updates FFXIV market-data columns (H-N)
# for specific Leve rows across all 8 job sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 309533.25
$ws.Cells.Item(86, 10).Value = 2000
$ws.Cells.Item(86, 12).Value = 2000
$ws.Cells.Item(86, 14).Value = -4246

$ws.Cells.Item(89, 8).Value = 309533.25
$ws.Cells.Item(89, 10).Value = 2000
$ws.Cells.Item(89, 12).Value = 10000
$ws.Cells.Item(89, 14).Value = -21232

$ws.Cells.Item(94, 8).Value = 3327.8333
$ws.Cells.Item(94, 9).Value = 2993.4
$ws.Cells.Item(94, 11).Value = 2993.4
$ws.Cells.Item(94, 13).Value = -2542.4

$ws.Cells.Item(97, 8).Value = 1400
$ws.Cells.Item(97, 10).Value = 1400
$ws.Cells.Item(97, 12).Value = 4200
$ws.Cells.Item(97, 14).Value = -5192

$ws.Cells.Item(111, 8).Value = 1167
$ws.Cells.Item(111, 9).Value = 367.6
$ws.Cells.Item(111, 10).Value = 3165.5
$ws.Cells.Item(111, 11).Value = 1102.8
$ws.Cells.Item(111, 12).Value = 9496.5
$ws.Cells.Item(111, 13).Value = 1964.2
$ws.Cells.Item(111, 14).Value = -15630.5

$ws.Cells.Item(117, 8).Value = 47742
$ws.Cells.Item(117, 10).Value = 47742
$ws.Cells.Item(117, 12).Value = 47742
$ws.Cells.Item(117, 14).Value = -56920

$ws.Cells.Item(129, 8).Value = 874
$ws.Cells.Item(129, 10).Value = 894.931
$ws.Cells.Item(129, 12).Value = 2684.793
$ws.Cells.Item(129, 14).Value = -12684.793

$ws.Cells.Item(132, 8).Value = 941.65216
$ws.Cells.Item(132, 9).Value = 876.8095
$ws.Cells.Item(132, 10).Value = 1622.5
$ws.Cells.Item(132, 11).Value = 2630.4285
$ws.Cells.Item(132, 12).Value = 4867.5
$ws.Cells.Item(132, 13).Value = -100.4285
$ws.Cells.Item(132, 14).Value = -9927.5

$ws.Cells.Item(137, 8).Value = 2245.4783
$ws.Cells.Item(137, 9).Value = 1328.3
$ws.Cells.Item(137, 10).Value = 2951
$ws.Cells.Item(137, 11).Value = 3984.9
$ws.Cells.Item(137, 12).Value = 8853
$ws.Cells.Item(137, 13).Value = -1434.9
$ws.Cells.Item(137, 14).Value = -13953

$ws.Cells.Item(138, 8).Value = 2925.2036
$ws.Cells.Item(138, 10).Value = 2866.5806
$ws.Cells.Item(138, 12).Value = 8599.7418
$ws.Cells.Item(138, 14).Value = -18879.7418

$ws.Cells.Item(141, 8).Value = 2728.9167
$ws.Cells.Item(141, 9).Value = 1194.2222
$ws.Cells.Item(141, 11).Value = 3582.6666
$ws.Cells.Item(141, 13).Value = 1597.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3785.675
$ws.Cells.Item(32, 9).Value = 3731.3428
$ws.Cells.Item(32, 10).Value = 4166
$ws.Cells.Item(32, 11).Value = 3731.3428
$ws.Cells.Item(32, 12).Value = 4166
$ws.Cells.Item(32, 13).Value = -3444.3428
$ws.Cells.Item(32, 14).Value = -4740

$ws.Cells.Item(61, 8).Value = 3541
$ws.Cells.Item(61, 9).Value = 1164
$ws.Cells.Item(61, 11).Value = 1164
$ws.Cells.Item(61, 13).Value = -952

$ws.Cells.Item(74, 8).Value = 1578.0526
$ws.Cells.Item(74, 9).Value = 1444.7
$ws.Cells.Item(74, 10).Value = 1726.2222
$ws.Cells.Item(74, 11).Value = 1444.7
$ws.Cells.Item(74, 12).Value = 1726.2222
$ws.Cells.Item(74, 13).Value = -570.7
$ws.Cells.Item(74, 14).Value = -3474.2222

$ws.Cells.Item(77, 8).Value = 1578.0526
$ws.Cells.Item(77, 9).Value = 1444.7
$ws.Cells.Item(77, 10).Value = 1726.2222
$ws.Cells.Item(77, 11).Value = 7223.5
$ws.Cells.Item(77, 12).Value = 8631.110999999999
$ws.Cells.Item(77, 13).Value = -2855.5
$ws.Cells.Item(77, 14).Value = -17367.111

$ws.Cells.Item(132, 8).Value = 1537.2273
$ws.Cells.Item(132, 9).Value = 953.5278
$ws.Cells.Item(132, 11).Value = 2860.5834
$ws.Cells.Item(132, 13).Value = -330.5834

$ws.Cells.Item(136, 8).Value = 3541
$ws.Cells.Item(136, 9).Value = 1164
$ws.Cells.Item(136, 11).Value = 3492
$ws.Cells.Item(136, 13).Value = -942

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2081.4482
$ws.Cells.Item(20, 9).Value = 2010.04
$ws.Cells.Item(20, 11).Value = 2010.04
$ws.Cells.Item(20, 13).Value = -1763.04

$ws.Cells.Item(82, 8).Value = 22564.25
$ws.Cells.Item(82, 9).Value = 13419
$ws.Cells.Item(82, 10).Value = 50000
$ws.Cells.Item(82, 11).Value = 13419
$ws.Cells.Item(82, 12).Value = 50000
$ws.Cells.Item(82, 13).Value = -13036
$ws.Cells.Item(82, 14).Value = -50766

$ws.Cells.Item(85, 8).Value = 22564.25
$ws.Cells.Item(85, 9).Value = 13419
$ws.Cells.Item(85, 10).Value = 50000
$ws.Cells.Item(85, 11).Value = 13419
$ws.Cells.Item(85, 12).Value = 50000
$ws.Cells.Item(85, 13).Value = -12093
$ws.Cells.Item(85, 14).Value = -52652

$ws.Cells.Item(97, 8).Value = 17042.2
$ws.Cells.Item(97, 9).Value = 5474.3335
$ws.Cells.Item(97, 10).Value = 21999.857
$ws.Cells.Item(97, 11).Value = 5474.3335
$ws.Cells.Item(97, 12).Value = 21999.857
$ws.Cells.Item(97, 13).Value = -4483.3335
$ws.Cells.Item(97, 14).Value = -23981.857

$ws.Cells.Item(99, 8).Value = 1916.75
$ws.Cells.Item(99, 9).Value = 1672.4166
$ws.Cells.Item(99, 10).Value = 2649.75
$ws.Cells.Item(99, 11).Value = 1672.4166
$ws.Cells.Item(99, 12).Value = 2649.75
$ws.Cells.Item(99, 13).Value = -174.4166
$ws.Cells.Item(99, 14).Value = -5645.75

$ws.Cells.Item(134, 8).Value = 5573.3945
$ws.Cells.Item(134, 9).Value = 6213.7
$ws.Cells.Item(134, 10).Value = 3172.25
$ws.Cells.Item(134, 11).Value = 18641.1
$ws.Cells.Item(134, 12).Value = 9516.75
$ws.Cells.Item(134, 13).Value = -16106.1
$ws.Cells.Item(134, 14).Value = -14586.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2428.4666
$ws.Cells.Item(31, 9).Value = 2283.3
$ws.Cells.Item(31, 10).Value = 2718.8
$ws.Cells.Item(31, 11).Value = 2283.3
$ws.Cells.Item(31, 12).Value = 2718.8
$ws.Cells.Item(31, 13).Value = -1988.3
$ws.Cells.Item(31, 14).Value = -3308.8

$ws.Cells.Item(34, 8).Value = 2428.4666
$ws.Cells.Item(34, 9).Value = 2283.3
$ws.Cells.Item(34, 10).Value = 2718.8
$ws.Cells.Item(34, 11).Value = 2283.3
$ws.Cells.Item(34, 12).Value = 2718.8
$ws.Cells.Item(34, 13).Value = -2081.3
$ws.Cells.Item(34, 14).Value = -3122.8

$ws.Cells.Item(99, 8).Value = 2182.3333
$ws.Cells.Item(99, 9).Value = 2132
$ws.Cells.Item(99, 11).Value = 2132
$ws.Cells.Item(99, 13).Value = -634

$ws.Cells.Item(126, 8).Value = 2182.3333
$ws.Cells.Item(126, 9).Value = 2132
$ws.Cells.Item(126, 11).Value = 6396
$ws.Cells.Item(126, 13).Value = -3926

$ws.Cells.Item(132, 8).Value = 1905.3667
$ws.Cells.Item(132, 9).Value = 1024.1666
$ws.Cells.Item(132, 10).Value = 5430.1665
$ws.Cells.Item(132, 11).Value = 3072.4998
$ws.Cells.Item(132, 12).Value = 16290.4995
$ws.Cells.Item(132, 13).Value = -542.4998000000001
$ws.Cells.Item(132, 14).Value = -21350.4995

$ws.Cells.Item(134, 8).Value = 2166.4138
$ws.Cells.Item(134, 9).Value = 1933.875
$ws.Cells.Item(134, 10).Value = 3282.6
$ws.Cells.Item(134, 11).Value = 5801.625
$ws.Cells.Item(134, 12).Value = 9847.799999999999
$ws.Cells.Item(134, 13).Value = -3266.625
$ws.Cells.Item(134, 14).Value = -14917.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 878.5
$ws.Cells.Item(122, 10).Value = 1148.7273
$ws.Cells.Item(122, 12).Value = 10338.5457
$ws.Cells.Item(122, 14).Value = -15238.5457

$ws.Cells.Item(131, 8).Value = 865.71
$ws.Cells.Item(131, 10).Value = 871.6667
$ws.Cells.Item(131, 12).Value = 2615.0001
$ws.Cells.Item(131, 14).Value = -12695.0001

$ws.Cells.Item(132, 8).Value = 990.0769
$ws.Cells.Item(132, 10).Value = 1082.8182
$ws.Cells.Item(132, 12).Value = 9745.363799999999
$ws.Cells.Item(132, 14).Value = -14805.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6485.5713
$ws.Cells.Item(70, 9).Value = 5599.75
$ws.Cells.Item(70, 10).Value = 7666.6665
$ws.Cells.Item(70, 11).Value = 5599.75
$ws.Cells.Item(70, 12).Value = 7666.6665
$ws.Cells.Item(70, 13).Value = -5329.75
$ws.Cells.Item(70, 14).Value = -8206.666499999999

$ws.Cells.Item(73, 8).Value = 6485.5713
$ws.Cells.Item(73, 9).Value = 5599.75
$ws.Cells.Item(73, 10).Value = 7666.6665
$ws.Cells.Item(73, 11).Value = 5599.75
$ws.Cells.Item(73, 12).Value = 7666.6665
$ws.Cells.Item(73, 13).Value = -4663.75
$ws.Cells.Item(73, 14).Value = -9538.666499999999

$ws.Cells.Item(97, 8).Value = 1875.5555
$ws.Cells.Item(97, 9).Value = 1875.5555
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 1875.5555
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -1379.5555
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 2573155.2
$ws.Cells.Item(126, 9).Value = 6175496.5
$ws.Cells.Item(126, 10).Value = 79226.92
$ws.Cells.Item(126, 11).Value = 18526489.5
$ws.Cells.Item(126, 12).Value = 237680.76
$ws.Cells.Item(126, 13).Value = -18524019.5
$ws.Cells.Item(126, 14).Value = -242620.76

$ws.Cells.Item(136, 8).Value = 10654.637
$ws.Cells.Item(136, 10).Value = 10654.637
$ws.Cells.Item(136, 12).Value = 31963.911
$ws.Cells.Item(136, 14).Value = -37063.911

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1298.8334
$ws.Cells.Item(100, 9).Value = 1298.8334
$ws.Cells.Item(100, 11).Value = 1298.8334
$ws.Cells.Item(100, 13).Value = -757.8334

$ws.Cells.Item(136, 8).Value = 2479
$ws.Cells.Item(136, 9).Value = 1341.7142
$ws.Cells.Item(136, 11).Value = 4025.1426
$ws.Cells.Item(136, 13).Value = -1475.1426

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 19158750
$ws.Cells.Item(136, 9).Value = 29241004
$ws.Cells.Item(136, 11).Value = 87723012
$ws.Cells.Item(136, 13).Value = -87720462
